$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$newRange = $lastPara.Range
$newRange.Collapse(0)
$newRange.InsertParagraphAfter()
$newRange.Collapse(0)
$newRange.InsertAfter("EN LOS FILTROS DE BUSQUEDA, CUANDO ESTA ULTIMA ES EXACTA, PERMITE FILTRAR POR CAMPOS NUMERICOS, (CASO CONTRARIO NO, NO PERMITIMOS HACER CONSULTAS LIKE CON CAMPOS NUMERICOS)")
